$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("index")
$ws3 = $wb.Worksheets.Item("problems_district")

# --- "index" sheet: sort the district rows (A2:J20) descending by column B,
# exactly like Data > Sort was used in Excel, leaving the "Total" summary row
# (originally row 21) out of the sorted range ---
$sortObj = $ws1.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws1.Range("B2"), $null, 2, $null, 2) | Out-Null
$sortObj.SetRange($ws1.Range("A2:J20"))
$sortObj.Header = 2
$sortObj.Apply()

# Move the "Total" row down one row, leaving row 21 blank as a gap above it
$ws1.Range("A21:J21").Cut($ws1.Range("A22:J22"))

# --- "problems_district" sheet: update its remembered selection and make it
# the non-active sheet ---
$ws3.Activate()
$ws3.Range("A5:XFD5,A11:XFD11").Select() | Out-Null

# --- "index" sheet becomes the active/selected sheet with a new selection ---
$ws1.Activate()
$ws1.Range("L25").Select() | Out-Null
